# Add data for 2021-12-05
# Updates the "through" date (sheet name + header label) from Nov 26 to
# Nov 27, and refreshes the carjacking counts that changed with the new
# data pull (current-month column B plus several backfilled historical
# months).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet and update the matching column-B header label.
$ws.Name = "Through 2021-11-27"
$ws.Range("B1").Value = "November 2021 (through November 27)"

# Cell -> new value updates (neighborhood/month intersections).
$updates = @{
    "M2"   = 19
    "X2"   = 6
    "AT2"  = 5
    "BE2"  = 2
    "M3"   = 16
    "M5"   = 6
    "B6"   = 12
    "B8"   = 7
    "B9"   = 2
    "M12"  = 8
    "M13"  = 5
    "M14"  = 2
    "B16"  = 4
    "AI31" = 1
    "B38"  = 2
    "AI38" = 1
    "M39"  = 2
    "BE48" = 5
    "M50"  = 1
    "BE51" = 2
    "M62"  = 2
    "BP91" = 1
    "B98"  = 3
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
